# Auto-generated edit script: swaps the full row contents for row pairs
# (25,26), (29,31), (30,32) while keeping row numbers fixed, matching the
# upstream diff exactly (re-ordered species records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# --- Row 25 ---
$ws.Range("A25").Value = 130887120
$ws.Range("B25").Value = 79243
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = 'Garnlav'
$ws.Range("G25").Value = 'Alectoria sarmentosa'
$ws.Range("H25").Value = '(Ach.) Ach.'
$ws.Range("I25").Value = ''
$ws.Range("J25").Value = ''
$ws.Range("K25").Value = 'med apothecier'
$ws.Range("L25").Value = ''
$ws.Range("M25").Value = ''
$ws.Range("N25").Value = ''
$ws.Range("Q25").Value = 495991
$ws.Range("R25").Value = 7016264
$ws.Range("AC25").Value = 'Långväxta fertila bålar med apothecier.'
$ws.Range("AF25").Value = ''
$ws.Range("AI25").Value = ''
$ws.Range("AJ25").Value = 'gran'
$ws.Range("AK25").Value = 'Picea abies'
$ws.Range("AM25").Value = 'Gren på levande träd'
$ws.Range("AO25").Value = 'Branch on living tree # Picea abies'

# --- Row 26 ---
$ws.Range("A26").Value = 130887107
$ws.Range("B26").Value = 57881
$ws.Range("E26").Value = 100049
$ws.Range("F26").Value = 'Spillkråka'
$ws.Range("G26").Value = 'Dryocopus martius'
$ws.Range("H26").Value = '(Linnaeus, 1758)'
$ws.Range("I26").Value = "'1"
$ws.Range("J26").Value = ''
$ws.Range("K26").Value = ''
$ws.Range("L26").Value = ''
$ws.Range("M26").Value = 'lockläte, övriga läten'
$ws.Range("N26").Value = 'observerad'
$ws.Range("Q26").Value = 495982
$ws.Range("R26").Value = 7016398
$ws.Range("AC26").Value = '1 eller möjligen 2 individer som sågs och hördes hacka i stående döda granar och uttryckte lockläte emellanåt. Hyfsat gott om stående död ved på/kring fyndplatsen.'
$ws.Range("AF26").Value = ''
$ws.Range("AI26").Value = 'Äldre bitvis flerskiktad granskog.'
$ws.Range("AJ26").Value = ''
$ws.Range("AK26").Value = ''
$ws.Range("AM26").Value = ''
$ws.Range("AO26").Value = ''

# --- Row 29 ---
$ws.Range("A29").Value = 130887096
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = 'Tretåig hackspett'
$ws.Range("G29").Value = 'Picoides tridactylus'
$ws.Range("H29").Value = '(Linnaeus, 1758)'
$ws.Range("J29").Value = ''
$ws.Range("L29").Value = ''
$ws.Range("M29").Value = 'äldre spår'
$ws.Range("Q29").Value = 495870
$ws.Range("R29").Value = 7016247
$ws.Range("AC29").Value = 'Ringhack, äldre, på en levande gran nära kant mot ungskog.'
$ws.Range("AF29").Value = ''
$ws.Range("AM29").Value = 'Trädstam på levande träd'
$ws.Range("AO29").Value = 'Stem on living tree # Picea abies'

# --- Row 30 ---
$ws.Range("A30").Value = 130887086
$ws.Range("B30").Value = 57884
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = 'Tretåig hackspett'
$ws.Range("G30").Value = 'Picoides tridactylus'
$ws.Range("H30").Value = '(Linnaeus, 1758)'
$ws.Range("J30").Value = ''
$ws.Range("L30").Value = ''
$ws.Range("M30").Value = 'färska spår'
$ws.Range("Q30").Value = 495683
$ws.Range("R30").Value = 7016064
$ws.Range("AC30").Value = 'Ringhack, färska, på en gran nära en mindre väg.'
$ws.Range("AF30").Value = ''
$ws.Range("AM30").Value = 'Trädstam på levande träd'
$ws.Range("AO30").Value = 'Stem on living tree # Picea abies'

# --- Row 31 ---
$ws.Range("A31").Value = 130887132
$ws.Range("B31").Value = 79243
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = 'Garnlav'
$ws.Range("G31").Value = 'Alectoria sarmentosa'
$ws.Range("H31").Value = '(Ach.) Ach.'
$ws.Range("J31").Value = ''
$ws.Range("L31").Value = ''
$ws.Range("M31").Value = ''
$ws.Range("Q31").Value = 496008
$ws.Range("R31").Value = 7016340
$ws.Range("AC31").Value = 'Långväxta bålar på gran. Ca 50 cm lång hängande bål.'
$ws.Range("AF31").Value = ''
$ws.Range("AM31").Value = 'Gren på levande träd'
$ws.Range("AO31").Value = 'Branch on living tree # Picea abies'

# --- Row 32 ---
$ws.Range("A32").Value = 130887127
$ws.Range("B32").Value = 79243
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = 'Garnlav'
$ws.Range("G32").Value = 'Alectoria sarmentosa'
$ws.Range("H32").Value = '(Ach.) Ach.'
$ws.Range("J32").Value = ''
$ws.Range("L32").Value = ''
$ws.Range("M32").Value = ''
$ws.Range("Q32").Value = 495831
$ws.Range("R32").Value = 7016258
$ws.Range("AC32").Value = 'På en gren av gran på ca 2 m höjd.'
$ws.Range("AF32").Value = ''
$ws.Range("AM32").Value = 'Gren på levande träd'
$ws.Range("AO32").Value = 'Branch on living tree # Picea abies'

Write-Host "Row swap edits applied."
